$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- was row 8's D/M/N/O/P/Q/S values
$ws.Range("D2").Value = 44400
$ws.Range("M2").Value = 100
$ws.Range("Q2").Value = "$/caja 14 kilos"

# Row 3 <- was row 4's D/M values
$ws.Range("D3").Value = 44309
$ws.Range("M3").Value = 300

# Row 4 <- was row 2's D/M/N/O/P/S values
$ws.Range("D4").Value = 44208
$ws.Range("M4").Value = 210
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("S4").Value = 714

# Row 6 <- was row 7's D/M/N/O/P/S values
$ws.Range("D6").Value = 44351
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("S6").Value = 714

# Row 7 <- was row 3's D/M/N/O/P/S values
$ws.Range("D7").Value = 44162
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 7000
$ws.Range("O7").Value = 7000
$ws.Range("P7").Value = 7000
$ws.Range("S7").Value = 500

# Row 8 <- was row 6's D/M/N/O/P/Q/S values
$ws.Range("D8").Value = 44176
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 7000
$ws.Range("O8").Value = 7000
$ws.Range("P8").Value = 7000
$ws.Range("Q8").Value = "$/caja 14 kilos empedrada"
$ws.Range("S8").Value = 500
